$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1124.25
$ws.Range("I12").Value = 1249
$ws.Range("K12").Value = 1249
$ws.Range("M12").Value = -1079
$ws.Range("H38").Value = 2170.4
$ws.Range("J38").Value = 3372.75
$ws.Range("L38").Value = 10118.25
$ws.Range("N38").Value = -10862.25
$ws.Range("H40").Value = 6242.3335
$ws.Range("I40").Value = 2655.5557
$ws.Range("J40").Value = 9829.111000000001
$ws.Range("K40").Value = 2655.5557
$ws.Range("L40").Value = 9829.111000000001
$ws.Range("M40").Value = -2480.5557
$ws.Range("N40").Value = -10179.111
$ws.Range("H76").Value = 8498.916999999999
$ws.Range("I76").Value = 6597.8
$ws.Range("K76").Value = 6597.8
$ws.Range("M76").Value = -6282.8
$ws.Range("H79").Value = 8498.916999999999
$ws.Range("I79").Value = 6597.8
$ws.Range("K79").Value = 6597.8
$ws.Range("M79").Value = -5505.8
$ws.Range("H82").Value = 11383.286
$ws.Range("I82").Value = 11383.286
$ws.Range("K82").Value = 34149.858
$ws.Range("M82").Value = -33743.858
$ws.Range("H85").Value = 11383.286
$ws.Range("I85").Value = 11383.286
$ws.Range("K85").Value = 34149.858
$ws.Range("M85").Value = -32745.858
$ws.Range("H99").Value = 4153.1
$ws.Range("I99").Value = 2630.125
$ws.Range("J99").Value = 10245
$ws.Range("K99").Value = 7890.375
$ws.Range("L99").Value = 30735
$ws.Range("M99").Value = -6392.375
$ws.Range("N99").Value = -33731
$ws.Range("H111").Value = 65974.31
$ws.Range("I111").Value = 1892.2307
$ws.Range("J111").Value = 343663.34
$ws.Range("K111").Value = 5676.6921
$ws.Range("L111").Value = 1030990.02
$ws.Range("M111").Value = -2609.6921
$ws.Range("N111").Value = -1037124.02
$ws.Range("H132").Value = 1708.575
$ws.Range("I132").Value = 1708.575
$ws.Range("K132").Value = 5125.725
$ws.Range("M132").Value = -2595.725
$ws.Range("H137").Value = 7085.636
$ws.Range("I137").Value = 4325.095
$ws.Range("J137").Value = 11916.583
$ws.Range("K137").Value = 12975.285
$ws.Range("L137").Value = 35749.749
$ws.Range("M137").Value = -10425.285
$ws.Range("N137").Value = -40849.749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2914.3333
$ws.Range("I45").Value = 1622
$ws.Range("J45").Value = 5499
$ws.Range("K45").Value = 1622
$ws.Range("L45").Value = 5499
$ws.Range("M45").Value = -1245
$ws.Range("N45").Value = -6253
$ws.Range("H74").Value = 1787.3489
$ws.Range("I74").Value = 629.0303
$ws.Range("J74").Value = 5609.8
$ws.Range("K74").Value = 629.0303
$ws.Range("L74").Value = 5609.8
$ws.Range("M74").Value = 244.9697
$ws.Range("N74").Value = -7357.8
$ws.Range("H77").Value = 1787.3489
$ws.Range("I77").Value = 629.0303
$ws.Range("J77").Value = 5609.8
$ws.Range("K77").Value = 3145.1515
$ws.Range("L77").Value = 28049
$ws.Range("M77").Value = 1222.8485
$ws.Range("N77").Value = -36785
$ws.Range("H82").Value = 128083.336
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 128083.336
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 128083.336
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -128805.336
$ws.Range("H85").Value = 128083.336
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 128083.336
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 128083.336
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -130579.336
$ws.Range("H132").Value = 2738.7856
$ws.Range("I132").Value = 1195.0454
$ws.Range("K132").Value = 3585.1362
$ws.Range("M132").Value = -1055.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4142.4443
$ws.Range("I134").Value = 1426.5
$ws.Range("J134").Value = 9574.333000000001
$ws.Range("K134").Value = 4279.5
$ws.Range("L134").Value = 28722.999
$ws.Range("M134").Value = -1744.5
$ws.Range("N134").Value = -33792.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5598.2197
$ws.Range("I31").Value = 2203.8572
$ws.Range("J31").Value = 12909.154
$ws.Range("K31").Value = 2203.8572
$ws.Range("L31").Value = 12909.154
$ws.Range("M31").Value = -1908.8572
$ws.Range("N31").Value = -13499.154
$ws.Range("H34").Value = 5598.2197
$ws.Range("I34").Value = 2203.8572
$ws.Range("J34").Value = 12909.154
$ws.Range("K34").Value = 2203.8572
$ws.Range("L34").Value = 12909.154
$ws.Range("M34").Value = -2001.8572
$ws.Range("N34").Value = -13313.154
$ws.Range("H99").Value = 4099.2
$ws.Range("I99").Value = 1536.1111
$ws.Range("J99").Value = 7943.8335
$ws.Range("K99").Value = 1536.1111
$ws.Range("L99").Value = 7943.8335
$ws.Range("M99").Value = -38.11110000000008
$ws.Range("N99").Value = -10939.8335
$ws.Range("H105").Value = 2772.6206
$ws.Range("I105").Value = 2500.4443
$ws.Range("J105").Value = 3218
$ws.Range("K105").Value = 2500.4443
$ws.Range("L105").Value = 3218
$ws.Range("M105").Value = -753.4443000000001
$ws.Range("N105").Value = -6712
$ws.Range("H126").Value = 4099.2
$ws.Range("I126").Value = 1536.1111
$ws.Range("J126").Value = 7943.8335
$ws.Range("K126").Value = 4608.3333
$ws.Range("L126").Value = 23831.5005
$ws.Range("M126").Value = -2138.3333
$ws.Range("N126").Value = -28771.5005
$ws.Range("H134").Value = 4713.533
$ws.Range("I134").Value = 3031.8635
$ws.Range("J134").Value = 9338.125
$ws.Range("K134").Value = 9095.5905
$ws.Range("L134").Value = 28014.375
$ws.Range("M134").Value = -6560.5905
$ws.Range("N134").Value = -33084.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 437.5
$ws.Range("I8").Value = 437.5
$ws.Range("K8").Value = 1312.5
$ws.Range("M8").Value = -1173.5
$ws.Range("H13").Value = 150
$ws.Range("I13").Value = 138
$ws.Range("K13").Value = 414
$ws.Range("M13").Value = -246
$ws.Range("H37").Value = 77000
$ws.Range("J37").Value = 77000
$ws.Range("L37").Value = 231000
$ws.Range("N37").Value = -231224
$ws.Range("H68").Value = 2357.9412
$ws.Range("I68").Value = 3766.3333
$ws.Range("J68").Value = 2056.1428
$ws.Range("K68").Value = 11298.9999
$ws.Range("L68").Value = 6168.428400000001
$ws.Range("M68").Value = -10487.9999
$ws.Range("N68").Value = -7790.428400000001
$ws.Range("H71").Value = 2357.9412
$ws.Range("I71").Value = 3766.3333
$ws.Range("J71").Value = 2056.1428
$ws.Range("K71").Value = 33896.9997
$ws.Range("L71").Value = 18505.2852
$ws.Range("M71").Value = -29840.9997
$ws.Range("N71").Value = -26617.2852
$ws.Range("H86").Value = 1072.8
$ws.Range("I86").Value = 1091
$ws.Range("K86").Value = 3273
$ws.Range("M86").Value = -2087
$ws.Range("H89").Value = 1072.8
$ws.Range("I89").Value = 1091
$ws.Range("K89").Value = 9819
$ws.Range("M89").Value = -3891
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496
$ws.Range("H107").Value = 1247.3334
$ws.Range("I107").Value = 1179.3334
$ws.Range("K107").Value = 3538.0002
$ws.Range("M107").Value = -1618.0002
$ws.Range("H131").Value = 1138095.2
$ws.Range("J131").Value = 1820563.8
$ws.Range("L131").Value = 5461691.4
$ws.Range("N131").Value = -5471771.4
$ws.Range("H134").Value = 1274.8889
$ws.Range("I134").Value = 772.5714
$ws.Range("J134").Value = 3033
$ws.Range("K134").Value = 2317.7142
$ws.Range("L134").Value = 9099
$ws.Range("M134").Value = 2752.2858
$ws.Range("N134").Value = -19239
$ws.Range("H137").Value = 2488.1177
$ws.Range("I137").Value = 1457.6
$ws.Range("J137").Value = 2917.5
$ws.Range("K137").Value = 4372.799999999999
$ws.Range("L137").Value = 8752.5
$ws.Range("M137").Value = 727.2000000000007
$ws.Range("N137").Value = -18952.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4415.0835
$ws.Range("I126").Value = 3310.125
$ws.Range("K126").Value = 9930.375
$ws.Range("M126").Value = -7460.375
$ws.Range("H132").Value = 2578.848
$ws.Range("I132").Value = 2068.805
$ws.Range("J132").Value = 6761.2
$ws.Range("K132").Value = 6206.414999999999
$ws.Range("L132").Value = 20283.6
$ws.Range("M132").Value = -3676.414999999999
$ws.Range("N132").Value = -25343.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2449.5
$ws.Range("I22").Value = 1975
$ws.Range("J22").Value = 2686.75
$ws.Range("K22").Value = 1975
$ws.Range("L22").Value = 2686.75
$ws.Range("M22").Value = -1680
$ws.Range("N22").Value = -3276.75
$ws.Range("H27").Value = 2449.5
$ws.Range("I27").Value = 1975
$ws.Range("J27").Value = 2686.75
$ws.Range("K27").Value = 1975
$ws.Range("L27").Value = 2686.75
$ws.Range("M27").Value = -1868
$ws.Range("N27").Value = -2900.75
$ws.Range("H40").Value = 5930.2085
$ws.Range("I40").Value = 4483.5
$ws.Range("K40").Value = 4483.5
$ws.Range("M40").Value = -4347.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3945.75
$ws.Range("I122").Value = 3622.6155
$ws.Range("K122").Value = 10867.8465
$ws.Range("M122").Value = -8417.8465

Write-Host "Applied all changes"